$p = $ppt.ActivePresentation

# Slide 2 (sldId 486): title "Session 1" -> "Session 2"
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 9).Text = "Session 2"

# Slide 3 (sldId 512): title "Session 1" -> "Session 2"
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 9).Text = "Session 2"

# Slide 19 (sldId 516): title "CSS" -> "CSS – updated to 3-tier format"
$s = $p.Slides.Item(19)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 3).Text = "CSS – updated to 3-tier format"
